# Atualização de bases das ligas, do dia: 10-06-2024 às 07:08
# Swap the data (columns B through AD) between pairs of rows whose
# "id" (column B) values were transposed in the source data, while
# leaving the row index in column A untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $rowA, $rowB, $firstCol, $lastCol) {
    $rangeA = $ws.Range($ws.Cells.Item($rowA, $firstCol), $ws.Cells.Item($rowA, $lastCol))
    $rangeB = $ws.Range($ws.Cells.Item($rowB, $firstCol), $ws.Cells.Item($rowB, $lastCol))
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

$firstCol = 2   # column B
$lastCol  = 30  # column AD

Swap-RowData $ws 41 42 $firstCol $lastCol
Swap-RowData $ws 47 48 $firstCol $lastCol
Swap-RowData $ws 137 138 $firstCol $lastCol
Swap-RowData $ws 155 156 $firstCol $lastCol
